$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change the Ledger Generation Date (B5): 2020-10-02 -> 2020-10-05 ---
# Force text formatting first so Excel doesn't auto-convert the date-like
# string into a date serial number, then drop back to the default/no style
# so the cell's style index matches what it was before (no "s" attribute).
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2020-10-05"
$ws.Range("B5").Style = "Normal"

# --- Insert a new data row above the TOTAL row (new row 14) ---
# This pushes the old row 14 (TOTAL) down to row 15 and auto-extends the
# sheet dimension to A1:E15.
$ws.Rows(14).Insert()

# Fill in the new row's values before formatting. A14 holds a date-looking
# string, so force text formatting first to avoid an auto date conversion.
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2020-10-01"
$ws.Range("B14").Value = "TEST"
$ws.Range("C14").Value = "CN-0005"
$ws.Range("D14").Value = "CREDITNOTE"
$ws.Range("E14").Value = 20

# Apply the same formatting as the row above (style s=2, thin border) onto
# the whole new row, matching the other CREDITNOTE detail rows.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122, -4142, $false, $false)

# --- Update the TOTAL formula so it includes the new row ---
$ws.Range("E15").Formula = "=SUM(E10:E14)"
